$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering of player rows (A2:C17) per the diff.
$data = @(
    ,@("CJ McCollum", "PG,SG", "New Orleans Pelicans")
    ,@("RJ Barrett", "SF,PF", "Toronto Raptors")
    ,@("Tobias Harris", "SF,PF", "Detroit Pistons")
    ,@("Brandon Boston Jr.", "SG,SF,PF", "New Orleans Pelicans")
    ,@("Jimmy Butler", "SF,PF", "Miami Heat")
    ,@("Keyonte George", "PG,SG", "Utah Jazz")
    ,@("Joel Embiid", "C", "Philadelphia 76ers")
    ,@("John Collins", "PF,C", "Utah Jazz")
    ,@("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder")
    ,@("Mike Conley", "PG", "Minnesota Timberwolves")
    ,@("Jordan Poole", "PG,SG", "Washington Wizards")
    ,@("Kyrie Irving", "PG,SG", "Dallas Mavericks")
    ,@("Shai Gilgeous-Alexander", "PG", "Oklahoma City Thunder")
    ,@("Zach LaVine", "SG,SF", "Chicago Bulls")
    ,@("Dennis Schröder", "PG", "Brooklyn Nets")
    ,@("Lauri Markkanen", "SF,PF", "Utah Jazz")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

Write-Output "done"